$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell A1 from "Name" to "COLLABORATORS"
$ws.Range("A1").Value = "COLLABORATORS"

# Update the active cell/selection to A5
$ws.Range("A5").Select()
